# Applies the FFXIV Leve-profit market-data refresh captured in the commit diff.
# All target cells are static cached numeric values (no formulas in the sheets),
# so each change is a direct Range.Value assignment; a few cells are removed
# outright (ClearContents) where the new data has no corresponding entry.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2845.25
$ws.Range("I12").Value = 3460
$ws.Range("J12").Value = 1001
$ws.Range("K12").Value = 3460
$ws.Range("L12").Value = 1001
$ws.Range("M12").Value = -3290
$ws.Range("N12").Value = -1341
$ws.Range("H62").Value = 2338.2354
$ws.Range("I62").Value = 2166.6667
$ws.Range("K62").Value = 2166.6667
$ws.Range("M62").Value = -1542.6667
$ws.Range("H65").Value = 2338.2354
$ws.Range("I65").Value = 2166.6667
$ws.Range("K65").Value = 10833.3335
$ws.Range("M65").Value = -7713.333500000001
$ws.Range("H86").Value = 8749
$ws.Range("I86").Value = 2299.7144
$ws.Range("J86").Value = 13765.111
$ws.Range("K86").Value = 2299.7144
$ws.Range("L86").Value = 13765.111
$ws.Range("M86").Value = -1176.7144
$ws.Range("N86").Value = -16011.111
$ws.Range("H89").Value = 8749
$ws.Range("I89").Value = 2299.7144
$ws.Range("J89").Value = 13765.111
$ws.Range("K89").Value = 11498.572
$ws.Range("L89").Value = 68825.55500000001
$ws.Range("M89").Value = -5882.572
$ws.Range("N89").Value = -80057.55500000001
$ws.Range("H92").Value = 1086
$ws.Range("I92").Value = 894.6
$ws.Range("K92").Value = 894.6
$ws.Range("M92").Value = 353.4
$ws.Range("H98").Value = 305.11765
$ws.Range("I98").Value = 305.11765
$ws.Range("K98").Value = 305.11765
$ws.Range("M98").Value = 1192.88235
$ws.Range("H122").Value = 305.11765
$ws.Range("I122").Value = 305.11765
$ws.Range("K122").Value = 915.3529500000001
$ws.Range("M122").Value = 1534.64705
$ws.Range("H129").Value = 839.0571
$ws.Range("J129").Value = 839.0571
$ws.Range("L129").Value = 2517.1713
$ws.Range("N129").Value = -12517.1713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10009
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("H20").Value = 10009
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H32").Value = 20016.771
$ws.Range("I32").Value = 20445.223
$ws.Range("J32").Value = 12304.667
$ws.Range("K32").Value = 20445.223
$ws.Range("L32").Value = 12304.667
$ws.Range("M32").Value = -20158.223
$ws.Range("N32").Value = -12878.667
$ws.Range("H37").Value = 29990
$ws.Range("J37").Value = 29990
$ws.Range("L37").Value = 29990
$ws.Range("N37").Value = -30536
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976
$ws.Range("H55").Value = 21220.6
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 24025.75
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 24025.75
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -24655.75
$ws.Range("H61").Value = 2529.1904
$ws.Range("I61").Value = 1950.7222
$ws.Range("K61").Value = 1950.7222
$ws.Range("M61").Value = -1738.7222
$ws.Range("H74").Value = 55559060
$ws.Range("I74").Value = 76926820
$ws.Range("K74").Value = 76926820
$ws.Range("M74").Value = -76925946
$ws.Range("H77").Value = 55559060
$ws.Range("I77").Value = 76926820
$ws.Range("K77").Value = 384634100
$ws.Range("M77").Value = -384629732
$ws.Range("H80").Value = 50845
$ws.Range("J80").Value = 50845
$ws.Range("L80").Value = 50845
$ws.Range("N80").Value = -52841
$ws.Range("H83").Value = 50845
$ws.Range("J83").Value = 50845
$ws.Range("L83").Value = 152535
$ws.Range("N83").Value = -162519
$ws.Range("H88").Value = 168913.67
$ws.Range("I88").Value = 2003
$ws.Range("K88").Value = 2003
$ws.Range("M88").Value = -1597
$ws.Range("H91").Value = 168913.67
$ws.Range("I91").Value = 2003
$ws.Range("K91").Value = 2003
$ws.Range("M91").Value = -599
$ws.Range("H102").Value = 1888.625
$ws.Range("I102").Value = 1703.3334
$ws.Range("K102").Value = 1703.3334
$ws.Range("M102").Value = -81.33339999999998
$ws.Range("H122").Value = 2061.3333
$ws.Range("I122").Value = 2031.5
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 6094.5
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -3644.5
$ws.Range("N122").Value = -11800
$ws.Range("H132").Value = 15242.784
$ws.Range("I132").Value = 1467.9615
$ws.Range("J132").Value = 47801.453
$ws.Range("K132").Value = 4403.8845
$ws.Range("L132").Value = 143404.359
$ws.Range("M132").Value = -1873.8845
$ws.Range("N132").Value = -148464.359
$ws.Range("H136").Value = 2529.1904
$ws.Range("I136").Value = 1950.7222
$ws.Range("K136").Value = 5852.1666
$ws.Range("M136").Value = -3302.1666
$ws.Range("M9").ClearContents()
$ws.Range("M20").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 73
$ws.Range("H86").Value = 1475.4375
$ws.Range("I86").Value = 1387.3948
$ws.Range("J86").Value = 1810
$ws.Range("K86").Value = 1387.3948
$ws.Range("L86").Value = 1810
$ws.Range("M86").Value = -264.3948
$ws.Range("N86").Value = -4056
$ws.Range("H89").Value = 1475.4375
$ws.Range("I89").Value = 1387.3948
$ws.Range("J89").Value = 1810
$ws.Range("K89").Value = 6936.974
$ws.Range("L89").Value = 9050
$ws.Range("M89").Value = -1320.974
$ws.Range("N89").Value = -20282
$ws.Range("H105").Value = 3335667.2
$ws.Range("I105").Value = 1867.4445
$ws.Range("K105").Value = 1867.4445
$ws.Range("M105").Value = -120.4445000000001
$ws.Range("H107").Value = 717.8
$ws.Range("I107").Value = 686.44446
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 686.44446
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1233.55554
$ws.Range("N107").Value = -4840
$ws.Range("N13").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H31").Value = 15704.869
$ws.Range("I31").Value = 26147.154
$ws.Range("J31").Value = 2129.9
$ws.Range("K31").Value = 26147.154
$ws.Range("L31").Value = 2129.9
$ws.Range("M31").Value = -25852.154
$ws.Range("N31").Value = -2719.9
$ws.Range("H34").Value = 15704.869
$ws.Range("I34").Value = 26147.154
$ws.Range("J34").Value = 2129.9
$ws.Range("K34").Value = 26147.154
$ws.Range("L34").Value = 2129.9
$ws.Range("M34").Value = -25945.154
$ws.Range("N34").Value = -2533.9
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 20258
$ws.Range("J107").Value = 330
$ws.Range("L107").Value = 990
$ws.Range("N107").Value = -4830
$ws.Range("H131").Value = 697.58
$ws.Range("I131").Value = 575
$ws.Range("J131").Value = 702.6875
$ws.Range("K131").Value = 1725
$ws.Range("L131").Value = 2108.0625
$ws.Range("M131").Value = 3315
$ws.Range("N131").Value = -12188.0625
$ws.Range("H139").Value = 1649.1666
$ws.Range("I139").Value = 1160.5
$ws.Range("K139").Value = 3481.5
$ws.Range("M139").Value = 1658.5
$ws.Range("H140").Value = 1747.3077
$ws.Range("I140").Value = 1470.909
$ws.Range("K140").Value = 4412.727000000001
$ws.Range("M140").Value = 767.2729999999992

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2985519.5
$ws.Range("I70").Value = 14690.7
$ws.Range("J70").Value = 5686273
$ws.Range("K70").Value = 14690.7
$ws.Range("L70").Value = 5686273
$ws.Range("M70").Value = -14420.7
$ws.Range("N70").Value = -5686813
$ws.Range("H73").Value = 2985519.5
$ws.Range("I73").Value = 14690.7
$ws.Range("J73").Value = 5686273
$ws.Range("K73").Value = 14690.7
$ws.Range("L73").Value = 5686273
$ws.Range("M73").Value = -13754.7
$ws.Range("N73").Value = -5688145
$ws.Range("H102").Value = 1014.5
$ws.Range("I102").Value = 1018
$ws.Range("J102").Value = 999.8
$ws.Range("K102").Value = 1018
$ws.Range("L102").Value = 999.8
$ws.Range("M102").Value = 604
$ws.Range("N102").Value = -4243.8
$ws.Range("H126").Value = 3790.7632
$ws.Range("I126").Value = 3006.4783
$ws.Range("K126").Value = 9019.4349
$ws.Range("M126").Value = -6549.4349

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3360.4443
$ws.Range("I22").Value = 3185.5
$ws.Range("J22").Value = 3500.4
$ws.Range("K22").Value = 3185.5
$ws.Range("L22").Value = 3500.4
$ws.Range("M22").Value = -2890.5
$ws.Range("N22").Value = -4090.4
$ws.Range("H27").Value = 3360.4443
$ws.Range("I27").Value = 3185.5
$ws.Range("J27").Value = 3500.4
$ws.Range("K27").Value = 3185.5
$ws.Range("L27").Value = 3500.4
$ws.Range("M27").Value = -3078.5
$ws.Range("N27").Value = -3714.4
$ws.Range("H100").Value = 1984.963
$ws.Range("I100").Value = 1789.8
$ws.Range("J100").Value = 2542.5715
$ws.Range("K100").Value = 1789.8
$ws.Range("L100").Value = 2542.5715
$ws.Range("M100").Value = -1248.8
$ws.Range("N100").Value = -3624.5715
$ws.Range("H132").Value = 1561.8529
$ws.Range("I132").Value = 974.13043
$ws.Range("K132").Value = 2922.39129
$ws.Range("M132").Value = -392.39129

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 40446
$ws.Range("J118").Value = 40446
$ws.Range("L118").Value = 40446
$ws.Range("N118").Value = -43760
